$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply each cell update as plain text, matching the original inline-string cell
# format (NumberFormat "@" forces text interpretation for numeric-looking values;
# ClearFormats afterwards removes the temporary style so cells keep their original
# (unstyled) appearance, matching the source workbook).

$r = $ws.Range("D2")
$r.NumberFormat = '@'
$r.Value = '61.909.30'
$r.ClearFormats()

$r = $ws.Range("E2")
$r.NumberFormat = '@'
$r.Value = '  -1.63%  '
$r.ClearFormats()

$r = $ws.Range("D3")
$r.NumberFormat = '@'
$r.Value = '3.413.73'
$r.ClearFormats()

$r = $ws.Range("E3")
$r.NumberFormat = '@'
$r.Value = '  -1.53%  '
$r.ClearFormats()

$r = $ws.Range("E4")
$r.NumberFormat = '@'
$r.Value = '  +0.00%  '
$r.ClearFormats()

$r = $ws.Range("D5")
$r.NumberFormat = '@'
$r.Value = '405.04'
$r.ClearFormats()

$r = $ws.Range("E5")
$r.NumberFormat = '@'
$r.Value = '  -0.92%  '
$r.ClearFormats()

$r = $ws.Range("D6")
$r.NumberFormat = '@'
$r.Value = '132.01'
$r.ClearFormats()

$r = $ws.Range("E6")
$r.NumberFormat = '@'
$r.Value = '  +0.26%  '
$r.ClearFormats()

$r = $ws.Range("D7")
$r.NumberFormat = '@'
$r.Value = '0.591'
$r.ClearFormats()

$r = $ws.Range("E7")
$r.NumberFormat = '@'
$r.Value = '  -1.62%  '
$r.ClearFormats()

$r = $ws.Range("D8")
$r.NumberFormat = '@'
$r.Value = '0.999'
$r.ClearFormats()

$r = $ws.Range("E8")
$r.NumberFormat = '@'
$r.Value = '  +0.06%  '
$r.ClearFormats()

$r = $ws.Range("D9")
$r.NumberFormat = '@'
$r.Value = '0.686'
$r.ClearFormats()

$r = $ws.Range("E9")
$r.NumberFormat = '@'
$r.Value = '  -1.48%  '
$r.ClearFormats()

$r = $ws.Range("E10")
$r.NumberFormat = '@'
$r.Value = '  -2.82%  '
$r.ClearFormats()

$r = $ws.Range("D11")
$r.NumberFormat = '@'
$r.Value = '41.83'
$r.ClearFormats()

$r = $ws.Range("E11")
$r.NumberFormat = '@'
$r.Value = '  -2.85%  '
$r.ClearFormats()

$r = $ws.Range("E12")
$r.NumberFormat = '@'
$r.Value = '  -0.96%  '
$r.ClearFormats()

$r = $ws.Range("D13")
$r.NumberFormat = '@'
$r.Value = '8.42'
$r.ClearFormats()

$r = $ws.Range("E13")
$r.NumberFormat = '@'
$r.Value = '  -4.45%  '
$r.ClearFormats()

$r = $ws.Range("D14")
$r.NumberFormat = '@'
$r.Value = '19.80'
$r.ClearFormats()

$r = $ws.Range("E14")
$r.NumberFormat = '@'
$r.Value = '  -1.66%  '
$r.ClearFormats()

$r = $ws.Range("D15")
$r.NumberFormat = '@'
$r.Value = '3.366.08'
$r.ClearFormats()

$r = $ws.Range("E15")
$r.NumberFormat = '@'
$r.Value = '  -2.55%  '
$r.ClearFormats()

$r = $ws.Range("D16")
$r.NumberFormat = '@'
$r.Value = '11.77'
$r.ClearFormats()

$r = $ws.Range("E16")
$r.NumberFormat = '@'
$r.Value = '  +8.51%  '
$r.ClearFormats()

$r = $ws.Range("D17")
$r.NumberFormat = '@'
$r.Value = '61.890.23'
$r.ClearFormats()

$r = $ws.Range("E17")
$r.NumberFormat = '@'
$r.Value = '  -1.55%  '
$r.ClearFormats()

$r = $ws.Range("E18")
$r.NumberFormat = '@'
$r.Value = '  -3.20%  '
$r.ClearFormats()

$r = $ws.Range("D19")
$r.NumberFormat = '@'
$r.Value = '0.0000142'
$r.ClearFormats()

$r = $ws.Range("E19")
$r.NumberFormat = '@'
$r.Value = '  +1.61%  '
$r.ClearFormats()

$r = $ws.Range("D20")
$r.NumberFormat = '@'
$r.Value = '3.17'
$r.ClearFormats()

$r = $ws.Range("E20")
$r.NumberFormat = '@'
$r.Value = '  -5.11%  '
$r.ClearFormats()

$r = $ws.Range("D21")
$r.NumberFormat = '@'
$r.Value = '83.53'
$r.ClearFormats()

$r = $ws.Range("E21")
$r.NumberFormat = '@'
$r.Value = '  +0.77%  '
$r.ClearFormats()

$r = $ws.Range("D22")
$r.NumberFormat = '@'
$r.Value = '311.55'
$r.ClearFormats()

$r = $ws.Range("E22")
$r.NumberFormat = '@'
$r.Value = '  -0.87%  '
$r.ClearFormats()

$r = $ws.Range("D23")
$r.NumberFormat = '@'
$r.Value = '12.83'
$r.ClearFormats()

$r = $ws.Range("E23")
$r.NumberFormat = '@'
$r.Value = '  -2.77%  '
$r.ClearFormats()

$r = $ws.Range("D24")
$r.NumberFormat = '@'
$r.Value = '3.14'
$r.ClearFormats()

$r = $ws.Range("E24")
$r.NumberFormat = '@'
$r.Value = '  -0.73%  '
$r.ClearFormats()

$r = $ws.Range("E25")
$r.NumberFormat = '@'
$r.Value = '  +9.99%  '
$r.ClearFormats()

$r = $ws.Range("D26")
$r.NumberFormat = '@'
$r.Value = '29.62'
$r.ClearFormats()

$r = $ws.Range("E26")
$r.NumberFormat = '@'
$r.Value = '  -2.85%  '
$r.ClearFormats()

$r = $ws.Range("D27")
$r.NumberFormat = '@'
$r.Value = '8.10'
$r.ClearFormats()

$r = $ws.Range("E27")
$r.NumberFormat = '@'
$r.Value = '  -1.87%  '
$r.ClearFormats()

$r = $ws.Range("D28")
$r.NumberFormat = '@'
$r.Value = '7.69'
$r.ClearFormats()

$r = $ws.Range("E28")
$r.NumberFormat = '@'
$r.Value = '  +0.72%  '
$r.ClearFormats()

$r = $ws.Range("D29")
$r.NumberFormat = '@'
$r.Value = '2.75'
$r.ClearFormats()

$r = $ws.Range("E29")
$r.NumberFormat = '@'
$r.Value = '  +5.72%  '
$r.ClearFormats()

$r = $ws.Range("E30")
$r.NumberFormat = '@'
$r.Value = '  -2.36%  '
$r.ClearFormats()

$r = $ws.Range("D31")
$r.NumberFormat = '@'
$r.Value = '0.115'
$r.ClearFormats()

$r = $ws.Range("E31")
$r.NumberFormat = '@'
$r.Value = '  -2.34%  '
$r.ClearFormats()

$r = $ws.Range("D32")
$r.NumberFormat = '@'
$r.Value = '42.68'
$r.ClearFormats()

$r = $ws.Range("E32")
$r.NumberFormat = '@'
$r.Value = '  -4.43%  '
$r.ClearFormats()

$r = $ws.Range("D34")
$r.NumberFormat = '@'
$r.Value = '11.34'
$r.ClearFormats()

$r = $ws.Range("E34")
$r.NumberFormat = '@'
$r.Value = '  -4.02%  '
$r.ClearFormats()

$r = $ws.Range("D35")
$r.NumberFormat = '@'
$r.Value = '0.0484'
$r.ClearFormats()

$r = $ws.Range("E35")
$r.NumberFormat = '@'
$r.Value = '  -1.95%  '
$r.ClearFormats()

$r = $ws.Range("D36")
$r.NumberFormat = '@'
$r.Value = '51.32'
$r.ClearFormats()

$r = $ws.Range("E36")
$r.NumberFormat = '@'
$r.Value = '  -2.44%  '
$r.ClearFormats()

$r = $ws.Range("E37")
$r.NumberFormat = '@'
$r.Value = '  -0.02%  '
$r.ClearFormats()

$r = $ws.Range("D38")
$r.NumberFormat = '@'
$r.Value = '0.327'
$r.ClearFormats()

$r = $ws.Range("E38")
$r.NumberFormat = '@'
$r.Value = '  +13.53%  '
$r.ClearFormats()

$r = $ws.Range("D39")
$r.NumberFormat = '@'
$r.Value = '3.38'
$r.ClearFormats()

$r = $ws.Range("E39")
$r.NumberFormat = '@'
$r.Value = '  -5.74%  '
$r.ClearFormats()

$r = $ws.Range("D40")
$r.NumberFormat = '@'
$r.Value = '2.93'
$r.ClearFormats()

$r = $ws.Range("E40")
$r.NumberFormat = '@'
$r.Value = '  -3.19%  '
$r.ClearFormats()

$r = $ws.Range("D41")
$r.NumberFormat = '@'
$r.Value = '139.23'
$r.ClearFormats()

$r = $ws.Range("E41")
$r.NumberFormat = '@'
$r.Value = '  +2.34%  '
$r.ClearFormats()

$r = $ws.Range("D42")
$r.NumberFormat = '@'
$r.Value = '0.125'
$r.ClearFormats()

$r = $ws.Range("E42")
$r.NumberFormat = '@'
$r.Value = '  -0.83%  '
$r.ClearFormats()

$r = $ws.Range("E43")
$r.NumberFormat = '@'
$r.Value = '  -0.87%  '
$r.ClearFormats()

$r = $ws.Range("D44")
$r.NumberFormat = '@'
$r.Value = '3.97'
$r.ClearFormats()

$r = $ws.Range("E44")
$r.NumberFormat = '@'
$r.Value = '  -0.29%  '
$r.ClearFormats()

$r = $ws.Range("D45")
$r.NumberFormat = '@'
$r.Value = '16.74'
$r.ClearFormats()

$r = $ws.Range("E45")
$r.NumberFormat = '@'
$r.Value = '  -4.26%  '
$r.ClearFormats()

$r = $ws.Range("E46")
$r.NumberFormat = '@'
$r.Value = '  -0.97%  '
$r.ClearFormats()

$r = $ws.Range("D47")
$r.NumberFormat = '@'
$r.Value = '21.18'
$r.ClearFormats()

$r = $ws.Range("E47")
$r.NumberFormat = '@'
$r.Value = '  -4.50%  '
$r.ClearFormats()

$r = $ws.Range("D48")
$r.NumberFormat = '@'
$r.Value = '2.108.69'
$r.ClearFormats()

$r = $ws.Range("E48")
$r.NumberFormat = '@'
$r.Value = '  -3.32%  '
$r.ClearFormats()

$r = $ws.Range("E49")
$r.NumberFormat = '@'
$r.Value = '  -3.25%  '
$r.ClearFormats()

$r = $ws.Range("B50")
$r.NumberFormat = '@'
$r.Value = 'Fetch.AI'
$r.ClearFormats()

$r = $ws.Range("C50")
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$r.ClearFormats()

$r = $ws.Range("D50")
$r.NumberFormat = '@'
$r.Value = '1.78'
$r.ClearFormats()

$r = $ws.Range("E50")
$r.NumberFormat = '@'
$r.Value = '  +21.39%  '
$r.ClearFormats()

$r = $ws.Range("B51")
$r.NumberFormat = '@'
$r.Value = 'ThetaToken'
$r.ClearFormats()

$r = $ws.Range("C51")
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$r.ClearFormats()

$r = $ws.Range("D51")
$r.NumberFormat = '@'
$r.Value = '1.96'
$r.ClearFormats()

$r = $ws.Range("E51")
$r.NumberFormat = '@'
$r.Value = '  +4.49%  '
$r.ClearFormats()
